$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2010
$ws.Range("B2").Value = 23257.86610669455
$ws.Range("C2").Value = 11418.82435341108
$ws.Range("D2").Value = 13560.43311809023

$ws.Range("A3").Value = 2011
$ws.Range("B3").Value = 25192.59314690959
$ws.Range("C3").Value = 12245.94221027538
$ws.Range("D3").Value = 14257.77027193148

$ws.Range("A4").Value = 2012
$ws.Range("B4").Value = 26979.84043437291
$ws.Range("C4").Value = 13211.21379819802
$ws.Range("D4").Value = 15697.30052998581

$ws.Range("A5").Value = 2013
$ws.Range("B5").Value = 27976.04665397691
$ws.Range("C5").Value = 13739.79027459498
$ws.Range("D5").Value = 16072.90969847004

$ws.Range("A6").Value = 2014
$ws.Range("B6").Value = 29912.95664259601
$ws.Range("C6").Value = 14745.8907374447
$ws.Range("D6").Value = 16785.95219769734

$ws.Range("A7").Value = 2015
$ws.Range("B7").Value = 31911.5976261855
$ws.Range("C7").Value = 15975.15151436092
$ws.Range("D7").Value = 17673.68733391674

$ws.Range("A8").Value = 2016
$ws.Range("B8").Value = 34224.49724617274
$ws.Range("C8").Value = 17607.09315882314
$ws.Range("D8").Value = 18604.68706838786

$ws.Range("A9").Value = 2017
$ws.Range("B9").Value = 35211.13679492509
$ws.Range("C9").Value = 18274.86689908277
$ws.Range("D9").Value = 19513.7812222432

$ws.Range("A10").Value = 2018
$ws.Range("B10").Value = 36646.18938775777
$ws.Range("C10").Value = 19083.90609078256
$ws.Range("D10").Value = 20594.91845447646

$ws.Range("A11").Value = 2019
$ws.Range("B11").Value = 37893.92756225695
$ws.Range("C11").Value = 19561.56719637316
$ws.Range("D11").Value = 20960.20827855906

$ws.Range("A12").Value = 2020
$ws.Range("B12").Value = 40040.12552822272
$ws.Range("C12").Value = 20906.9003924803
$ws.Range("D12").Value = 21326.36830874367

$ws.Range("A13").Value = 2021
$ws.Range("B13").Value = 44932.83670743866
$ws.Range("C13").Value = 22962.50837336785
$ws.Range("D13").Value = 23155.15872464016

$ws.Range("A14").Value = 2022
$ws.Range("B14").Value = 48783.70604975082
$ws.Range("C14").Value = 24759.83037637985
$ws.Range("D14").Value = 25299.05850862464

$ws.Range("A15").Value = 2023
$ws.Range("B15").Value = 51300.70579350938
$ws.Range("C15").Value = 26237.41536180414
$ws.Range("D15").Value = 26006.98661973922
